$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure re-tagging after re-curation
$ws.Range("B2").Value = "iaest-measure:vehiculos-en-el-hogar"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: medida/dim swap between columns B and D
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "dim"

# Row 4: type/URI columns re-curated
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"

# Row 5 (mapping file references) no longer needed - remove entirely
$ws.Range("B5:E5").Clear()
